$wb = $excel.ActiveWorkbook

# Sheet "Arkusz6" (6th worksheet) has its "time" column (B2:B5) converted
# from text values like "15 min" to plain numbers (15, 47, 120, 373).
$ws6 = $wb.Worksheets.Item(6)
$ws6.Activate()
$ws6.Range("B2").Value = 15
$ws6.Range("B3").Value = 47
$ws6.Range("B4").Value = 120
$ws6.Range("B5").Value = 373
$ws6.Range("B4").Select()

# Arkusz3 remains the active/selected tab; update its remembered
# selection to B16.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("B16").Select()
